$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 34.852685675233843
$ws.Range("C2").Value = 16.082241287026495
$ws.Range("D2").Value = 0.46143477828035678
$ws.Range("E2").Value = 30.762583279862866
$ws.Range("F2").Value = 15.359820453489798
$ws.Range("G2").Value = 0.49930203565005249
$ws.Range("H2").Value = 339
$ws.Range("I2").Value = 303.5

$ws.Range("B3").Value = 35.056145633056971
$ws.Range("C3").Value = 16.252334883549981
$ws.Range("D3").Value = 0.46360872223854754
$ws.Range("E3").Value = 30.589320366890298
$ws.Range("F3").Value = 15.21192787004836
$ws.Range("G3").Value = 0.49729538569656034
$ws.Range("H3").Value = 339
$ws.Range("I3").Value = 303

$ws.Range("B4").Value = 34.758320304648343
$ws.Range("C4").Value = 15.978548903702665
$ws.Range("D4").Value = 0.45970428845969868
$ws.Range("E4").Value = 30.74272736618687
$ws.Range("F4").Value = 15.359180370335979
$ws.Range("G4").Value = 0.49960370097902063
$ws.Range("H4").Value = 338.5
$ws.Range("I4").Value = 303

$ws.Range("B5").Value = 34.89133841210343
$ws.Range("C5").Value = 16.09135863563845
$ws.Range("D5").Value = 0.46118490628197084
$ws.Range("E5").Value = 30.862960336249095
$ws.Range("F5").Value = 15.47113461294332
$ws.Range("G5").Value = 0.50128485551569713
$ws.Range("H5").Value = 338.5
$ws.Range("I5").Value = 303

$ws.Range("B6").Value = 35.020544640300422
$ws.Range("C6").Value = 16.204328356616671
$ws.Range("D6").Value = 0.46270920464124632
$ws.Range("E6").Value = 30.610743259516603
$ws.Range("F6").Value = 15.252512662761756
$ws.Range("G6").Value = 0.49827318903861922
$ws.Range("H6").Value = 338.5
$ws.Range("I6").Value = 302

$ws.Range("B7").Value = 35.11391883754348
$ws.Range("C7").Value = 16.278583163972016
$ws.Range("D7").Value = 0.46359346102284388
$ws.Range("E7").Value = 30.70026999890041
$ws.Range("F7").Value = 15.337890741797686
$ws.Range("G7").Value = 0.49960116775347718
$ws.Range("H7").Value = 338.5
$ws.Range("I7").Value = 302

$ws.Range("B8").Value = 34.720323321491925
$ws.Range("C8").Value = 15.933081566429149
$ws.Range("D8").Value = 0.45889784547503193
$ws.Range("E8").Value = 30.770627147941887
$ws.Range("F8").Value = 15.373363899942538
$ws.Range("G8").Value = 0.4996116532181501
$ws.Range("H8").Value = 338.5
$ws.Range("I8").Value = 302

$ws.Range("B9").Value = 34.782452001312727
$ws.Range("C9").Value = 15.985983610313303
$ws.Range("D9").Value = 0.45959909927310399
$ws.Range("E9").Value = 30.833659773795375
$ws.Range("F9").Value = 15.431302458765821
$ws.Range("G9").Value = 0.50046937573983463
$ws.Range("H9").Value = 338.5
$ws.Range("I9").Value = 302

$ws.Range("B10").Value = 34.848722231711442
$ws.Range("C10").Value = 16.040471173122416
$ws.Range("D10").Value = 0.46028864606479025
$ws.Range("E10").Value = 30.88621676200659
$ws.Range("F10").Value = 15.485614300039105
$ws.Range("G10").Value = 0.50137620995680177
$ws.Range("H10").Value = 338.5
$ws.Range("I10").Value = 302

$ws.Range("B11").Value = 34.902943652677102
$ws.Range("C11").Value = 16.086617324795494
$ws.Range("D11").Value = 0.46089571942341401
$ws.Range("E11").Value = 30.945583760443384
$ws.Range("F11").Value = 15.53188159928167
$ws.Range("G11").Value = 0.50190947178496959
$ws.Range("H11").Value = 338.5
$ws.Range("I11").Value = 302
